# Update "paises.xlsx" (sheet "Pais") with refreshed COVID country stats
# and a handful of label-order swaps among tied/adjacent countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "last updated" timestamp note in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Junio de 2020 a las 18:02"

# --- 2. Swap the ordering of a few tied/adjacent countries (column A labels) ---
# Moldavia / Ghana  (rows 57-58)
$ws.Range("A57").Value = "Ghana"
$ws.Range("A58").Value = "Moldavia"

# Tanzania / Reunion (rows 153-154)
$ws.Range("A153").Value = "Reunion"
$ws.Range("A154").Value = "Tanzania"

# Dominica / Fiyi (rows 202-203) - tied values, label order only
$ws.Range("A202").Value = "Fiyi"
$ws.Range("A203").Value = "Dominica"

# Islas Malvinas / Groenlandia (rows 208-209) - tied values, label order only
$ws.Range("A208").Value = "Groenlandia"
$ws.Range("A209").Value = "Islas Malvinas"

# Montserrat / Seychelles (rows 211-212)
$ws.Range("A211").Value = "Seychelles"
$ws.Range("A212").Value = "Montserrat"

# --- 3. Refresh the numeric columns (B:H) for the affected rows ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 2474928
$ws.Range("C4").Value = 12374
$ws.Range("D4").Value = 1040711
$ws.Range("E4").Value = 1309803
$ws.Range("G4").Value = 133
$ws.Range("H4").Value = 124414

# Row 7 - India
$ws.Range("B7").Value = 481179
$ws.Range("C7").Value = 8194
$ws.Range("D7").Value = 277765
$ws.Range("E7").Value = 188372
$ws.Range("G7").Value = 135
$ws.Range("H7").Value = 15042

# Row 22 - Canada
$ws.Range("B22").Value = 102573
$ws.Range("C22").Value = 331
$ws.Range("D22").Value = 65361
$ws.Range("E22").Value = 28711
$ws.Range("G22").Value = 17
$ws.Range("H22").Value = 8501

# Row 26 - Suecia
$ws.Range("B26").Value = 63890
$ws.Range("C26").Value = 230
$ws.Range("G26").Value = 21
$ws.Range("H26").Value = 5230

# Row 36 - Singapur
$ws.Range("D36").Value = 36604
$ws.Range("E36").Value = 6106

# Row 45 - Republica Dominicana
$ws.Range("B45").Value = 29141
$ws.Range("C45").Value = 510
$ws.Range("D45").Value = 16223
$ws.Range("E45").Value = 12220
$ws.Range("G45").Value = 7
$ws.Range("H45").Value = 698

# Row 50 - Barein
$ws.Range("E50").Value = 5523
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 70

# Row 54 - Kazajistan
$ws.Range("D54").Value = 12220
$ws.Range("E54").Value = 6925
$ws.Range("G54").Value = 4
$ws.Range("H54").Value = 140

# Row 57 - now Ghana (new, updated figures)
$ws.Range("B57").Value = 15473
$ws.Range("C57").Value = 460
$ws.Range("D57").Value = 11431
$ws.Range("E57").Value = 3947
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 95

# Row 58 - now Moldavia (carried-over figures)
$ws.Range("B58").Value = 15453
$ws.Range("C58").Value = 375
$ws.Range("D58").Value = 8599
$ws.Range("E58").Value = 6352
$ws.Range("G58").Value = 7
$ws.Range("H58").Value = 502

# Row 69 - Chequia
$ws.Range("B69").Value = 10830
$ws.Range("C69").Value = 53
$ws.Range("D69").Value = 7649
$ws.Range("E69").Value = 2836
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 345

# Row 79 - Tayikistan
$ws.Range("B79").Value = 5691
$ws.Range("C79").Value = 61
$ws.Range("D79").Value = 4267
$ws.Range("E79").Value = 1372

# Row 90 - Luxemburgo
$ws.Range("B90").Value = 4151
$ws.Range("C90").Value = 11
$ws.Range("D90").Value = 3968
$ws.Range("E90").Value = 73

# Row 95 - Grecia
$ws.Range("B95").Value = 3321
$ws.Range("C95").Value = 11
$ws.Range("E95").Value = 1756
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 191

# Row 98 - Somalia
$ws.Range("B98").Value = 2878
$ws.Range("C98").Value = 43
$ws.Range("D98").Value = 868
$ws.Range("E98").Value = 1920

# Row 105 - Albania
$ws.Range("E105").Value = 893
$ws.Range("G105").Value = 2
$ws.Range("H105").Value = 49

# Row 128 - Jordania
$ws.Range("B128").Value = 1086
$ws.Range("C128").Value = 15
$ws.Range("D128").Value = 797

# Row 132 - Cabo Verde
$ws.Range("B132").Value = 1003
$ws.Range("C132").Value = 4
$ws.Range("E132").Value = 433

# Row 153 - now Reunion (new, updated figures)
$ws.Range("B153").Value = 516
$ws.Range("C153").Value = 8
$ws.Range("D153").Value = 460
$ws.Range("E153").Value = 54
$ws.Range("G153").Value = 1
$ws.Range("H153").Value = 2

# Row 154 - now Tanzania (carried-over figures)
$ws.Range("B154").Value = 509
$ws.Range("D154").Value = 183
$ws.Range("E154").Value = 305
$ws.Range("H154").Value = 21

# Row 211 - now Seychelles
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

# Row 212 - now Montserrat
$ws.Range("D212").Value = 10
$ws.Range("H212").Value = 1
